# Update the "Förändrad" (Changed) date column (column C) for all data rows
# from 45202 (2023-10-03) to 45203 (2023-10-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 261
$col = 3  # Column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
